# Automatic update of files.
# - Bump the "Förändrad" (column C) date for every data row to 2023-09-17 (serial 45186).
# - Add the record's "Beteckning" (column A) as the 2nd HYPERLINK() argument for every
#   non-empty hyperlink formula cell in columns S, T, V, W, X, Y (friendly link text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1

# Column C = 3 ("Förändrad" date). Row 1 is the header row, data starts at row 2.
$dateCol = 3
$newDateSerial = 45186

# Hyperlink-bearing columns, keyed by their column number.
$linkCols = @(19, 20, 22, 23, 24, 25)   # S, T, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {

    $idCell = $ws.Cells.Item($r, 1)
    if ($idCell.Value2 -eq $null -or $idCell.Value2 -eq "") {
        continue
    }

    # --- update the "Förändrad" date ---
    $dateCell = $ws.Cells.Item($r, $dateCol)
    if ($dateCell.Value2 -ne $null -and $dateCell.Value2 -ne "") {
        $dateCell.Value = $newDateSerial
    }

    # --- add the link text (2nd HYPERLINK argument) ---
    $beteckning = $idCell.Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if (-not $cell.HasFormula) {
            continue
        }

        $formula = $cell.Formula
        if ($formula -notmatch "HYPERLINK\(") {
            continue
        }
        if ($formula -match ",") {
            # already has a second argument
            continue
        }

        $newFormula = $formula -replace '\)\s*$', (', "' + $beteckning + '")')
        $cell.Formula = $newFormula
    }
}
